# Add files via upload
#
# Row 18 (account 729833) was missing its Rep (col B) and Manager (col D)
# values. Fill them in so the row matches the others for that account.
#
# Write D18 before B18: new shared-string entries are appended in the
# order they are first written, and the target file expects
# "Fred Anderson" to land at shared-string index 29 (written/used first)
# and "Koepp Ltd" at index 30 (written second) -- i.e. D18 then B18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "Fred Anderson"
$ws.Range("B18").Value = "Koepp Ltd"

# Keep the selection where the author left it.
[void]$ws.Range("B18").Select()

# Best-effort: reflect the saved window position (cosmetic bookViews
# xWindow/yWindow on the workbook). Harmless no-op if unsupported.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 360
    $win.Top = 460
} catch {
}
